$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Activate()

# Clear the extra header cells D3:G3 (Success, Failed, Error, Payload)
$ws.Range("D3:G3").Clear()

# Move the selection back to A3 to match the saved view state
$ws.Range("A3").Select()
